$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.NumberFormat = "General"
    $r.Style = "Normal"
}

Set-TextValue "D2" '61.289.29'
Set-TextValue "E2" '  +7.60%  '
Set-TextValue "D3" '3.340.51'
Set-TextValue "E3" '  +2.46%  '
Set-TextValue "E4" '  -0.04%  '
Set-TextValue "D5" '411.67'
Set-TextValue "E5" '  +3.96%  '
Set-TextValue "D6" '116.11'
Set-TextValue "E6" '  +6.81%  '
Set-TextValue "D7" '3.334.43'
Set-TextValue "E7" '  +2.38%  '
Set-TextValue "E8" '  -2.24%  '
Set-TextValue "E9" '  -0.01%  '
Set-TextValue "D10" '0.628'
Set-TextValue "E10" '  +0.13%  '
Set-TextValue "E11" '  +18.34%  '
Set-TextValue "D12" '40.09'
Set-TextValue "E12" '  +2.01%  '
Set-TextValue "E13" '  -0.74%  '
Set-TextValue "D14" '3.866.80'
Set-TextValue "E14" '  +2.31%  '
Set-TextValue "D15" '8.32'
Set-TextValue "E15" '  -0.66%  '
Set-TextValue "D16" '19.14'
Set-TextValue "E16" '  -0.19%  '
Set-TextValue "D17" '3.326.68'
Set-TextValue "E17" '  +1.85%  '
Set-TextValue "D18" '61.084.97'
Set-TextValue "E18" '  +7.38%  '
Set-TextValue "E19" '  -2.34%  '
Set-TextValue "D20" '10.85'
Set-TextValue "E20" '  +0.69%  '
Set-TextValue "D21" '0.0000116'
Set-TextValue "E21" '  +6.56%  '
Set-TextValue "E22" '  +0.37%  '
Set-TextValue "D23" '12.55'
Set-TextValue "E23" '  -3.90%  '
Set-TextValue "D24" '294.67'
Set-TextValue "E24" '  -0.55%  '
Set-TextValue "D25" '74.25'
Set-TextValue "E25" '  -0.11%  '
Set-TextValue "E26" '  -1.86%  '
Set-TextValue "D27" '29.15'
Set-TextValue "E27" '  +3.37%  '
Set-TextValue "D28" '7.83'
Set-TextValue "E28" '  +7.73%  '
Set-TextValue "E29" '  -2.49%  '
Set-TextValue "E30" '  +2.14%  '
Set-TextValue "D31" '7.56'
Set-TextValue "E31" '  -1.42%  '
Set-TextValue "E32" '  +5.00%  '
Set-TextValue "D33" '42.69'
Set-TextValue "E33" '  +6.76%  '
Set-TextValue "D34" '0.999'
Set-TextValue "E34" '  -0.04%  '
Set-TextValue "B35" 'Cosmos'
Set-TextValue "C35" 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue "D35" '11.34'
Set-TextValue "E35" '  +0.70%  '
Set-TextValue "B36" 'Toncoin'
Set-TextValue "C36" 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue "D36" '2.53'
Set-TextValue "E36" '  +18.55%  '
Set-TextValue "E37" '  -0.84%  '
Set-TextValue "D38" '52.42'
Set-TextValue "E38" '  +1.72%  '
Set-TextValue "E39" '  -0.05%  '
Set-TextValue "D40" '3.08'
Set-TextValue "E40" '  +4.29%  '
Set-TextValue "D41" '3.45'
Set-TextValue "E41" '  -1.38%  '
Set-TextValue "D42" '135.03'
Set-TextValue "E42" '  -3.37%  '
Set-TextValue "E43" '  -1.37%  '
Set-TextValue "D44" '0.289'
Set-TextValue "E44" '  +2.32%  '
Set-TextValue "E45" '  +0.02%  '
Set-TextValue "D46" '3.87'
Set-TextValue "E46" '  -3.20%  '
Set-TextValue "D47" '16.39'
Set-TextValue "E47" '  -4.56%  '
Set-TextValue "D48" '2.24'
Set-TextValue "E48" '  +4.61%  '
Set-TextValue "D49" '21.18'
Set-TextValue "E49" '  -5.18%  '
Set-TextValue "D50" '2.153.83'
Set-TextValue "E50" '  -0.67%  '
Set-TextValue "D51" '3.671.11'
Set-TextValue "E51" '  +2.42%  '
